$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Astar")
$ws2 = $wb.Worksheets.Item("avg")

# ---------------------------------------------------------------------
# Sheet1 (Astar): fill in previously-empty D/E/F columns for several rows
# ---------------------------------------------------------------------

# Row 14-18 (Mini CS / Mini Physics block under Astar)
$ws1.Range("D14").Value = 10021
$ws1.Range("E14").Value = 85.05
$ws1.Range("F14").Value = 14

$ws1.Range("D15").Value = 192605
$ws1.Range("E15").Value = 80.23
$ws1.Range("F15").Value = 55

$ws1.Range("D16").Value = 3327
$ws1.Range("E16").Value = 79.33
$ws1.Range("F16").Value = 1

$ws1.Range("D17").Value = 363
$ws1.Range("E17").Value = 77.900000000000006
$ws1.Range("F17").Value = 2

$ws1.Range("D18").Value = 206
$ws1.Range("E18").Value = 77.900000000000006
$ws1.Range("F18").Value = 0.02

$ws1.Range("D19").Value = "-"
$ws1.Range("E19").Value = "-"
$ws1.Range("F19").Value = "-"

# Row 35-40 (DFS block)
$ws1.Range("D35").Value = 64
$ws1.Range("E35").Value = 75.56
$ws1.Range("F35").Value = 0.001

$ws1.Range("D36").Value = 10532
$ws1.Range("E36").Value = 76.13
$ws1.Range("F36").Value = 0.8

$ws1.Range("D37").Value = 163
$ws1.Range("E37").Value = 74.180000000000007
$ws1.Range("F37").Value = 0.01

$ws1.Range("D38").Value = 921
$ws1.Range("E38").Value = 70.31
$ws1.Range("F38").Value = 0.09

$ws1.Range("D39").Value = 10466
$ws1.Range("E39").Value = 75.41
$ws1.Range("F39").Value = 0.5

$ws1.Range("D40").Value = "-"
$ws1.Range("E40").Value = "-"
$ws1.Range("F40").Value = "-"

# Row 41-49 (UCS block)
$ws1.Range("D41").Value = 1260027
$ws1.Range("E41").Value = 84.46
$ws1.Range("F41").Value = 150

$ws1.Range("D42").Value = 135912
$ws1.Range("E42").Value = 84.46
$ws1.Range("F42").Value = 13

$ws1.Range("D43").Value = 21960
$ws1.Range("E43").Value = 83.01
$ws1.Range("F43").Value = 2

$ws1.Range("D44").Value = 1700926
$ws1.Range("E44").Value = 85.05
$ws1.Range("F44").Value = 300

$ws1.Range("D45").Value = 1134098
$ws1.Range("E45").Value = 80.23
$ws1.Range("F45").Value = 100

$ws1.Range("D46").Value = 52509
$ws1.Range("E46").Value = 79.33
$ws1.Range("F46").Value = 10

$ws1.Range("D47").Value = 957805
$ws1.Range("E47").Value = 77.930000000000007
$ws1.Range("F47").Value = 280

$ws1.Range("D48").Value = 957805
$ws1.Range("E48").Value = 77.930000000000007
$ws1.Range("F48").Value = 120

$ws1.Range("D49").Value = "-"
$ws1.Range("E49").Value = "-"
$ws1.Range("F49").Value = "-"

# Rows 53/54: swap the Algorithm label (B) and Avg value (E)
$ws1.Range("B53").Value = "Physics"
$ws1.Range("E53").Value = 78.817999999999998

$ws1.Range("B54").Value = "Mini CS"
$ws1.Range("E54").Value = 84.606999999999999

# New rows 55 and 56 (Upper Bound / Mini Math, Mini Physics)
$ws1.Range("A55").Value = "Upper Bound"
$ws1.Range("B55").Value = "Mini Math"
$ws1.Range("C55").Value = "-"
$ws1.Range("D55").Value = "-"
$ws1.Range("E55").Value = 86.566000000000003
$ws1.Range("F55").Value = "-"

$ws1.Range("A56").Value = "Upper Bound"
$ws1.Range("B56").Value = "Mini Physics"
$ws1.Range("C56").Value = "-"
$ws1.Range("D56").Value = "-"
$ws1.Range("E56").Value = 78.126999999999995
$ws1.Range("F56").Value = "-"

# Page setup: A4 paper (xlPaperA4 = 9), portrait orientation (xlPortrait = 1)
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# Sheet view scroll/selection state
$ws1.Application.ActiveWindow.ScrollRow = 34
$ws1.Application.ActiveWindow.ScrollColumn = 1
$ws1.Range("D48").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet2 (avg): only view-state (scroll position/selection) changes
# ---------------------------------------------------------------------
$ws2.Select() | Out-Null
$ws2.Application.ActiveWindow.ScrollRow = 28
$ws2.Application.ActiveWindow.ScrollColumn = 1
$ws2.Range("B47").Select() | Out-Null

# Reselect sheet1 as the active/visible tab (tabSelected=1 on Astar)
$ws1.Select() | Out-Null
